$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Chloroplèthe avec valeurs extrêmes dans pop-up" -> split into two
#    runs, with a proofErr spellStart/spellEnd pair wrapping the first
#    word ("Chloroplèthe"), matching the look of a Word spell-check
#    annotation around an unrecognised word.
# ---------------------------------------------------------------------
$target = $d.Content
$found = $target.Find.Execute("Chloroplèthe avec valeurs extrêmes dans pop-up")
if (-not $found) {
    throw "Could not find the Chloroplèthe cell text"
}
Write-Host "Found target text:" $found

$newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
          '<w:body>' +
            '<w:p w14:paraId="2ECC1CA8" w14:textId="77777777" w:rsidR="00492E45" w:rsidRDefault="00492E45" w:rsidP="00492E45">' +
              '<w:proofErr w:type="spellStart"/>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>' +
                  '<w:color w:val="000000"/>' +
                '</w:rPr>' +
                '<w:t>Chloroplèthe</w:t>' +
              '</w:r>' +
              '<w:proofErr w:type="spellEnd"/>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>' +
                  '<w:color w:val="000000"/>' +
                '</w:rPr>' +
                '<w:t xml:space="preserve"> avec valeurs extrêmes dans pop-up</w:t>' +
              '</w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

[void]$target.InsertXML($newXml)
Write-Host "Inserted spell-check-annotated runs for Chloroplèthe"

# ---------------------------------------------------------------------
# 2) Drop the two trailing blank rows of the table (after the
#    "Densité de population" / "Antoine" row).
# ---------------------------------------------------------------------
$table = $d.Tables.Item(1)
while ($table.Rows.Count -gt 6) {
    $table.Rows.Item($table.Rows.Count).Delete()
}
Write-Host "Table now has" $table.Rows.Count "rows"
